$d = $word.ActiveDocument

# --- Step 1: remove the "Yo q se" text from the last paragraph ---
# (this paragraph currently also hosts the "_GoBack" bookmark)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$textRange = $lastPara.Range
$textRange.MoveEnd(1, -1)   # wdCharacter = 1; exclude the paragraph mark
$textRange.Text = ""

# --- Step 2: move the "_GoBack" bookmark to the empty paragraph right ---
# --- after the first heading ("-ZER DA CMS? (IGOR)")                 ---
$targetPara = $d.Paragraphs(2)
$bmRange = $targetPara.Range
$bmRange.MoveEnd(1, -1)     # collapse onto the (empty) paragraph content

# Adding a bookmark with a name that already exists relocates it, since
# bookmark names must be unique within a document.
$d.Bookmarks.Add("_GoBack", $bmRange)
